$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Report ")

# Update the Unit Test "actual test cases" figures (Test/UnitTest update).
$ws.Range("E8").Value = 459
$ws.Range("C17").Value = 459
$ws.Range("G17").Value = 459

# Recalculate the dependent totals (E12, C20, G20, ...).
$excel.Calculate()

# Reflect the cursor position left by the edit.
$ws.Range("G18").Select()
